# Update column G ("K" - strikeouts) values for rows 2-32 on Sheet1.
# This corresponds to regenerating save_data with K (strikeouts) instead of
# Strike# as described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @(9,2,5,6,1,5,5,4,4,5,4,9,3,4,6,2,2,8,4,2,7,5,5,2,3,2,2,2,3,1,1)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newValues[$i]
}
